$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 73254
$ws.Range("B2").Value = "Felipe Marques"
$ws.Range("C2").Value = "Marketing"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 45103
$ws.Range("G2").Value = 8754.459999999999

# Row 3
$ws.Range("A3").Value = 27711
$ws.Range("B3").Value = "Srta. Rafaela Ramos"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("D3").Value = "Doenca"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45084
$ws.Range("G3").Value = 3431.49

# Row 4
$ws.Range("A4").Value = 15566
$ws.Range("B4").Value = "Clarice Novais"
$ws.Range("C4").Value = "Marketing"
$ws.Range("D4").Value = "Doenca"
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 45086
$ws.Range("G4").Value = 9548.889999999999

# Row 5
$ws.Range("A5").Value = 25354
$ws.Range("B5").Value = "Dr. Benício Dias"
$ws.Range("D5").Value = "Doenca"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 45100
$ws.Range("G5").Value = 7593.38

# Row 6
$ws.Range("A6").Value = 49242
$ws.Range("B6").Value = "Miguel Souza"
$ws.Range("D6").Value = "Viagem de negocios"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 45084
$ws.Range("G6").Value = 4536.76

# Row 7
$ws.Range("A7").Value = 50111
$ws.Range("B7").Value = "Sarah Alves"
$ws.Range("C7").Value = "Vendas"
$ws.Range("D7").Value = "Consulta medica"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 45104
$ws.Range("G7").Value = 7967.36

# Row 8
$ws.Range("A8").Value = 15166
$ws.Range("B8").Value = "Larissa Nogueira"
$ws.Range("D8").Value = "Viagem de negocios"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45094
$ws.Range("G8").Value = 6927.9

# Row 9
$ws.Range("A9").Value = 2208
$ws.Range("B9").Value = "Bella Pinto"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Doenca"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45081
$ws.Range("G9").Value = 5044.21

# Row 10
$ws.Range("A10").Value = 62820
$ws.Range("B10").Value = "Eduarda Câmara"
$ws.Range("C10").Value = "Engenharia"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45081
$ws.Range("G10").Value = 9375.73

# Row 11
$ws.Range("A11").Value = 63170
$ws.Range("B11").Value = "Dra. Ana Cecília Albuquerque"
$ws.Range("C11").Value = "Juridico"
$ws.Range("D11").Value = "Doenca"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 45088
$ws.Range("G11").Value = 5628.34
